$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 794, pushing the existing rows 794-823 down to 796-825.
$ws.Range("A794:T795").EntireRow.Insert()

# Row 794: new "Especial" quality record, 25/04/2023 (serial 45041)
$ws.Cells.Item(794, 1).Value = 9
$ws.Cells.Item(794, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(794, 3).Value = "Metropolitana"
$ws.Cells.Item(794, 4).Value = 45041
$ws.Cells.Item(794, 5).Value = 13
$ws.Cells.Item(794, 6).Value = "Fruta"
$ws.Cells.Item(794, 7).Value = 100101
$ws.Cells.Item(794, 8).Value = "Berries"
$ws.Cells.Item(794, 9).Value = 100101007
$ws.Cells.Item(794, 10).Value = "Kiwi"
$ws.Cells.Item(794, 11).Value = "Hayward"
$ws.Cells.Item(794, 12).Value = "Especial"
$ws.Cells.Item(794, 13).Value = 250
$ws.Cells.Item(794, 14).Value = 12000
$ws.Cells.Item(794, 15).Value = 12000
$ws.Cells.Item(794, 16).Value = 12000
$ws.Cells.Item(794, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(794, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(794, 19).Value = 1200
$ws.Cells.Item(794, 20).Value = 10

# Row 795: new "Primera" quality record, 25/04/2023 (serial 45041)
$ws.Cells.Item(795, 1).Value = 9
$ws.Cells.Item(795, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(795, 3).Value = "Metropolitana"
$ws.Cells.Item(795, 4).Value = 45041
$ws.Cells.Item(795, 5).Value = 13
$ws.Cells.Item(795, 6).Value = "Fruta"
$ws.Cells.Item(795, 7).Value = 100101
$ws.Cells.Item(795, 8).Value = "Berries"
$ws.Cells.Item(795, 9).Value = 100101007
$ws.Cells.Item(795, 10).Value = "Kiwi"
$ws.Cells.Item(795, 11).Value = "Hayward"
$ws.Cells.Item(795, 12).Value = "Primera"
$ws.Cells.Item(795, 13).Value = 280
$ws.Cells.Item(795, 14).Value = 10000
$ws.Cells.Item(795, 15).Value = 10000
$ws.Cells.Item(795, 16).Value = 10000
$ws.Cells.Item(795, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(795, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(795, 19).Value = 1000
$ws.Cells.Item(795, 20).Value = 10
